$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.379.76'
$ws.Range("D3").Value = '1.581.59'
$ws.Range("E3").Value = '  -0.11%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''213.01'
$ws.Range("E5").Value = '  +0.34%  '
$ws.Range("D6").Value = '''0.491'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("D8").Value = '''44.73'
$ws.Range("E8").Value = '  -4.71%  '
$ws.Range("D9").Value = '''23.87'
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("D12").Value = '''0.0895'
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '1.807.44'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '1.575.03'
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '''0.517'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").Value = '28.411.70'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '''62.00'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = '''229.53'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("D20").Value = '''7.45'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '0.0₃0689'
$ws.Range("E21").Value = '  -1.89%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = '''3.92'
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("D24").Value = '''9.12'
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("D26").Value = '''151.37'
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("D27").Value = '''15.03'
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("D28").Value = '''6.41'
$ws.Range("E28").Value = '  -1.60%  '
$ws.Range("E29").Value = '  -1.09%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").Value = '''0.0480'
$ws.Range("E31").Value = '  +2.68%  '
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("E33").Value = '  -0.92%  '
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("D35").Value = '1.400.23'
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("E36").Value = '  +7.54%  '
$ws.Range("E37").Value = '  -4.29%  '
$ws.Range("D38").Value = '''2.37'
$ws.Range("E38").Value = '  +0.44%  '
$ws.Range("D39").Value = '''2.66'
$ws.Range("E39").Value = '  +1.89%  '
$ws.Range("E40").Value = '  -1.03%  '
$ws.Range("D41").Value = '''0.521'
$ws.Range("E41").Value = '  -2.69%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '''0.788'
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").Value = '''0.0461'
$ws.Range("E45").Value = '  -2.58%  '
$ws.Range("D46").Value = '''5.44'
$ws.Range("E46").Value = '  -2.69%  '
$ws.Range("D47").Value = '''62.70'
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E48").Value = '  -5.68%  '
$ws.Range("D49").Value = '1.720.14'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '''86.31'
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -1.04%  '
